$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vessels")

$ws.Range("D2").Value = 166.67
$ws.Range("J2").Value = 166.67
$ws.Range("D3").Value = 333.33
$ws.Range("J3").Value = 166.67
$ws.Range("D4").Value = 541.67
$ws.Range("J4").Value = 83.33

$ws.Activate()
$ws.Range("D5").Select()
